# Updates the cryptos list values (Price and Volume(1h) columns) to match
# the latest scrape. A handful of rows (Maker/Cosmos/dogwifhat and
# ThetaToken/USDe) were re-ranked, so their Coin/Link/Price/Volume moved
# to a different row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Some Price values look like plain numbers (e.g. "599.57"). Force those
# particular cells to keep a Text format so Excel stores them as strings
# (matching the original inline-string cells) instead of converting them
# to numeric values.
$textAddresses = @("D5", "D6", "D10", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D28", "D32", "D36", "D37", "D39", "D40", "D41", "D43", "D44", "D45", "D47", "D48", "D49", "D51")
foreach ($addr in $textAddresses) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.842.28'
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Value = '3.294.97'
$ws.Range("E3").Value = '  +5.43%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '599.57'
$ws.Range("E5").Value = '  +1.34%  '
$ws.Range("D6").Value = '143.84'
$ws.Range("E6").Value = '  +6.53%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.289.69'
$ws.Range("E8").Value = '  +5.55%  '
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("D10").Value = '0.150'
$ws.Range("E10").Value = '  +3.05%  '
$ws.Range("D11").Value = '5.45'
$ws.Range("E11").Value = '  +2.70%  '
$ws.Range("D12").Value = '0.473'
$ws.Range("E12").Value = '  +3.67%  '
$ws.Range("D13").Value = '0.0000249'
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").Value = '35.02'
$ws.Range("E14").Value = '  +3.46%  '
$ws.Range("D15").Value = '3.825.37'
$ws.Range("E15").Value = '  +5.28%  '
$ws.Range("E16").Value = '  +1.35%  '
$ws.Range("D17").Value = '3.284.78'
$ws.Range("E17").Value = '  +5.28%  '
$ws.Range("D18").Value = '63.903.51'
$ws.Range("E18").Value = '  +1.34%  '
$ws.Range("D19").Value = '6.92'
$ws.Range("E19").Value = '  +3.76%  '
$ws.Range("D20").Value = '483.08'
$ws.Range("E20").Value = '  +2.18%  '
$ws.Range("D21").Value = '14.32'
$ws.Range("E21").Value = '  +1.56%  '
$ws.Range("D22").Value = '0.745'
$ws.Range("E22").Value = '  +7.59%  '
$ws.Range("D23").Value = '8.06'
$ws.Range("E23").Value = '  +6.03%  '
$ws.Range("D24").Value = '13.57'
$ws.Range("E24").Value = '  +4.87%  '
$ws.Range("D25").Value = '84.77'
$ws.Range("E25").Value = '  -2.55%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  +2.55%  '
$ws.Range("D28").Value = '7.31'
$ws.Range("E28").Value = '  +2.89%  '
$ws.Range("E29").Value = '  +4.33%  '
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("E31").Value = '  +6.92%  '
$ws.Range("D32").Value = '28.22'
$ws.Range("E32").Value = '  +4.33%  '
$ws.Range("E33").Value = '  +1.24%  '
$ws.Range("E34").Value = '  +2.12%  '
$ws.Range("E35").Value = '  +2.41%  '
$ws.Range("D36").Value = '6.03'
$ws.Range("E36").Value = '  +3.62%  '
$ws.Range("D37").Value = '53.13'
$ws.Range("E37").Value = '  +2.20%  '
$ws.Range("D38").Value = '0.0₃0742'
$ws.Range("E38").Value = '  +4.64%  '
$ws.Range("D39").Value = '0.0399'
$ws.Range("E39").Value = '  +3.05%  '
$ws.Range("D40").Value = '428.64'
$ws.Range("E40").Value = '  +2.54%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").Value = '2.81'
$ws.Range("E41").Value = '  +4.57%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '3.021.95'
$ws.Range("E42").Value = '  +5.67%  '
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").Value = '8.49'
$ws.Range("E43").Value = '  +3.60%  '
$ws.Range("D44").Value = '0.111'
$ws.Range("E44").Value = '  -4.45%  '
$ws.Range("D45").Value = '0.271'
$ws.Range("E45").Value = '  +6.16%  '
$ws.Range("E46").Value = '  +7.88%  '
$ws.Range("D47").Value = '26.36'
$ws.Range("E47").Value = '  +4.60%  '
$ws.Range("B48").Value = 'ThetaToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D48").Value = '2.36'
$ws.Range("E48").Value = '  +4.25%  '
$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").Value = '0.999'
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("E50").Value = '  +1.90%  '
$ws.Range("D51").Value = '123.37'
$ws.Range("E51").Value = '  +4.17%  '
